$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 (item 11): "بلاستر مترسيلك 2 سم" -> "بلاستر مترسيلك 2.5 سم", transactions "21:0" -> "27:0"
$ws.Range("C17").Value = "بلاستر مترسيلك 2.5 سم"
$ws.Range("H17").Value = "27:0"

# Row 18 (item 12, فرشاة اطفال ريتش ديلي): price "15.00"/"15.0000" -> "25.00"/"25.0000"
$ws.Range("N18").Value = "25.00"
# P18 is formatted with a numeric (0.00) number format, so a plain .Value assignment of a
# numeric-looking string would be auto-converted to a real number. Force text entry, then
# restore the original numeric display format so the cell keeps its original style/format.
$ws.Range("P18").NumberFormat = "@"
$ws.Range("P18").Value = "25.0000"
$ws.Range("P18").NumberFormat = "0.00"

# Total row reflects the price increase for row 18 (+10.00)
$ws.Range("P20").Value = 667.015

# Footer timestamp updated from 11:21 AM to 11:23 AM
$ws.Range("A21").Value = "Saturday, 27 September, 2025 11:23 AM"
